$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Relocate the signature block (old rows 28-29) down to new rows 53-54,
#    since the worker table below is growing and will occupy rows 16-48.
# ---------------------------------------------------------------------------
$ws.Range("B28:C28").Copy($ws.Range("B53:C53"))
$ws.Range("H28:J28").Copy($ws.Range("H53:J53"))
$ws.Range("B29:C29").Copy($ws.Range("B54:C54"))
$ws.Range("H29:J29").Copy($ws.Range("H54:J54"))

# Clear the old signature rows' contents (will be fully overwritten below anyway).
$ws.Range("B28:J29").ClearContents()

# ---------------------------------------------------------------------------
# 2. Propagate the two row styles used by the worker table:
#    - "normal" row style (as used by old rows 16-22)
#    - "final" row style (as used by old row 23 - heavier bottom border)
#    The new table keeps the same "final" style only on its very last row (48).
# ---------------------------------------------------------------------------
$ws.Range("B23:J23").Copy($ws.Range("B48:J48"))
$ws.Range("B16:J16").Copy($ws.Range("B17:J47"))

# ---------------------------------------------------------------------------
# 3. Fill in the worker / period detail rows (16-48).
#    5 workers, most with 6 months of arrears (2502-2507), the last one
#    with 9 (2411,2412,2501-2507).
# ---------------------------------------------------------------------------
$workers = @(
    @{ Doc = "1050958934"; Name = "LUIS ENRIQUE GALVIS FERNANDEZ"; Periods = @("2507","2506","2505","2504","2503","2502") },
    @{ Doc = "1047461100"; Name = "BORIS ENRIQUE ALVIS LOPEZ"; Periods = @("2507","2506","2505","2504","2503","2502") },
    @{ Doc = "1043636844"; Name = "MAURO JAVIER MOGOLLON SALAS"; Periods = @("2507","2506","2505","2504","2503","2502") },
    @{ Doc = "1091668249"; Name = "LUIS FERNANDO TARAZONA QUINTERO"; Periods = @("2507","2506","2505","2504","2503","2502") },
    @{ Doc = "1050963173"; Name = "DANIEL ANDRES PEREZ GOMEZ"; Periods = @("2507","2506","2505","2504","2503","2502","2501","2412","2411") }
)

$row = 16
foreach ($worker in $workers) {
    foreach ($period in $worker.Periods) {
        $ws.Cells.Item($row, 2).Value2 = "CC"
        $ws.Cells.Item($row, 3).Value2 = $worker.Doc
        $ws.Cells.Item($row, 4).Value2 = $worker.Name
        $ws.Cells.Item($row, 5).Value2 = $period
        $ws.Cells.Item($row, 6).Value2 = 52000
        $ws.Cells.Item($row, 7).Value2 = 1300000
        $row = $row + 1
    }
}

# ---------------------------------------------------------------------------
# 4. Update the summary header cells.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value2 = 1716000   # VALOR MORA total
$ws.Range("C13").Value2 = 5         # Cant. Trabajadores
$ws.Range("F13").Value2 = 9         # Cant. Periodos

# ---------------------------------------------------------------------------
# 5. Fix up merged cells for the relocated signature block.
# ---------------------------------------------------------------------------
$ws.Range("B53:C53").Merge()
$ws.Range("H53:J53").Merge()
$ws.Range("B54:C54").Merge()
$ws.Range("H54:J54").Merge()


